$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I ("target") held "deuteron" for every data row (2-10); update it to "d".
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 9).Value = "d"
}

# Header row (row 1) becomes bold + centered.
$header = $ws.Range("A1:K1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108

# Move the active selection to H15.
$ws.Range("H15").Select()
